$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove all existing data (coordinates table + comment strings) so the
# workbook starts from a clean sheet before writing the new Apple-logo
# coordinate data.
$ws.Cells.ClearContents() | Out-Null

# Re-write the header row.
$ws.Range("A1").Value = "X"
$ws.Range("B1").Value = "Y"
$ws.Range("C1").Value = "Radius"
$ws.Range("D1").Value = "Comment"

# New coordinate data (X, Y, Radius) - no comments this time.
$data = @(
    @(156, -239, 8),
    @(312, -239, 8),
    @(238, -98, 8),
    @(235, -83, 8),
    @(320, -157, 8),
    @(418, -287, 8),
    @(299, -396, 8),
    @(238, -514, 8),
    @(238, -396, 1),
    @(177, -396, 5),
    @(299, -396, 5),
    @(114, -344, 3),
    @(363, -344, 3),
    @(204, -277, 20),
    @(271, -277, 20),
    @(238, -298, 13)
)

$row = 2
foreach ($d in $data) {
    $ws.Cells.Item($row, 1).Value = $d[0]
    $ws.Cells.Item($row, 2).Value = $d[1]
    $ws.Cells.Item($row, 3).Value = $d[2]
    $row = $row + 1
}

# Narrow column D down to fit the new (shorter) content instead of the wide
# 24.44-character columns used for the old comment text.
$ws.Columns("D").ColumnWidth = 10.5

# Match the saved selection/active cell.
$ws.Range("D2").Select() | Out-Null
